$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-06-11 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-06-12 Wednesday", 2)

# Update the division problems in the table, cell by cell (row, column),
# since several old values repeat with different replacements depending
# on position.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="868÷5="},
    @{Row=1;  Col=2; Text="694÷9="},
    @{Row=1;  Col=3; Text="117÷6="},
    @{Row=1;  Col=4; Text="646÷3="},
    @{Row=1;  Col=5; Text="726÷3="},

    @{Row=5;  Col=1; Text="392÷8="},
    @{Row=5;  Col=2; Text="196÷7="},
    @{Row=5;  Col=3; Text="941÷3="},
    @{Row=5;  Col=4; Text="975÷9="},
    @{Row=5;  Col=5; Text="565÷3="},

    @{Row=9;  Col=1; Text="585÷6="},
    @{Row=9;  Col=2; Text="952÷8="},
    @{Row=9;  Col=3; Text="421÷7="},
    @{Row=9;  Col=4; Text="534÷6="},
    @{Row=9;  Col=5; Text="171÷8="},

    @{Row=13; Col=1; Text="731÷3="},
    @{Row=13; Col=2; Text="773÷8="},
    @{Row=13; Col=3; Text="939÷7="},
    @{Row=13; Col=4; Text="920÷2="},
    @{Row=13; Col=5; Text="375÷9="},

    @{Row=17; Col=1; Text="508÷2="},
    @{Row=17; Col=2; Text="865÷8="},
    @{Row=17; Col=3; Text="887÷4="},
    @{Row=17; Col=4; Text="293÷7="},
    @{Row=17; Col=5; Text="154÷7="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
